$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 11.188041
$ws.Cells.Item(2, 8).Value = 33.564123
$ws.Cells.Item(2, 9).Value = 0.1395103797998223
$ws.Cells.Item(2, 10).Value = 0.1395103797998223
$ws.Cells.Item(2, 13).Value = 4.114675333333333
$ws.Cells.Item(2, 14).Value = 12.344026
$ws.Cells.Item(2, 15).Value = 0.04794018332925969
$ws.Cells.Item(2, 16).Value = 0.04794018332925969
$ws.Cells.Item(2, 17).Value = 46.035156331022
$ws.Cells.Item(2, 18).Value = 414.316406979198
$ws.Cells.Item(2, 19).Value = 0.006688153183938127
$ws.Cells.Item(2, 20).Value = 0.006688153183938129

$ws.Cells.Item(3, 7).Value = 11.188041
$ws.Cells.Item(3, 8).Value = 33.564123
$ws.Cells.Item(3, 9).Value = 0.1395103797998223
$ws.Cells.Item(3, 10).Value = 0.1395103797998223
$ws.Cells.Item(3, 13).Value = 59.62659933333333
$ws.Cells.Item(3, 15).Value = 0.6947109727426806
$ws.Cells.Item(3, 16).Value = 0.6947109727426806
$ws.Cells.Item(3, 17).Value = 667.104838031906
$ws.Cells.Item(3, 18).Value = 6003.943542287154
$ws.Cells.Item(3, 19).Value = 0.09691939165843536
$ws.Cells.Item(3, 20).Value = 0.09691939165843537

$ws.Cells.Item(4, 7).Value = 11.188041
$ws.Cells.Item(4, 8).Value = 33.564123
$ws.Cells.Item(4, 9).Value = 0.1395103797998223
$ws.Cells.Item(4, 10).Value = 0.1395103797998223
$ws.Cells.Item(4, 13).Value = 21.65107466666666
$ws.Cells.Item(4, 14).Value = 64.95322399999999
$ws.Cells.Item(4, 15).Value = 0.252257202503176
$ws.Cells.Item(4, 16).Value = 0.252257202503176
$ws.Cells.Item(4, 17).Value = 242.233111064728
$ws.Cells.Item(4, 18).Value = 2180.097999582552
$ws.Cells.Item(4, 19).Value = 0.03519249812845877
$ws.Cells.Item(4, 20).Value = 0.03519249812845877

$ws.Cells.Item(5, 7).Value = 11.188041
$ws.Cells.Item(5, 8).Value = 33.564123
$ws.Cells.Item(5, 9).Value = 0.1395103797998223
$ws.Cells.Item(5, 10).Value = 0.1395103797998223
$ws.Cells.Item(5, 13).Value = 0.4370123333333333
$ws.Cells.Item(5, 14).Value = 1.311037
$ws.Cells.Item(5, 15).Value = 0.005091641424883797
$ws.Cells.Item(5, 16).Value = 0.005091641424883797
$ws.Cells.Item(5, 17).Value = 4.889311902839
$ws.Cells.Item(5, 18).Value = 44.00380712555101
$ws.Cells.Item(5, 19).Value = 0.0007103368289900467
$ws.Cells.Item(5, 20).Value = 0.0007103368289900469

$ws.Cells.Item(6, 9).Value = 0.4168441980730721
$ws.Cells.Item(6, 10).Value = 0.4168441980730722
$ws.Cells.Item(6, 13).Value = 4.114675333333333
$ws.Cells.Item(6, 14).Value = 12.344026
$ws.Cells.Item(6, 15).Value = 0.04794018332925969
$ws.Cells.Item(6, 16).Value = 0.04794018332925969
$ws.Cells.Item(6, 17).Value = 137.5488178837129
$ws.Cells.Item(6, 18).Value = 1237.939360953416
$ws.Cells.Item(6, 19).Value = 0.01998358727536131
$ws.Cells.Item(6, 20).Value = 0.01998358727536132

$ws.Cells.Item(7, 9).Value = 0.4168441980730721
$ws.Cells.Item(7, 10).Value = 0.4168441980730722
$ws.Cells.Item(7, 13).Value = 59.62659933333333
$ws.Cells.Item(7, 15).Value = 0.6947109727426806
$ws.Cells.Item(7, 16).Value = 0.6947109727426806
$ws.Cells.Item(7, 17).Value = 1993.247969355974
$ws.Cells.Item(7, 19).Value = 0.2895862383254865
$ws.Cells.Item(7, 20).Value = 0.2895862383254866

$ws.Cells.Item(8, 9).Value = 0.4168441980730721
$ws.Cells.Item(8, 10).Value = 0.4168441980730722
$ws.Cells.Item(8, 13).Value = 21.65107466666666
$ws.Cells.Item(8, 14).Value = 64.95322399999999
$ws.Cells.Item(8, 15).Value = 0.252257202503176
$ws.Cells.Item(8, 16).Value = 0.252257202503176
$ws.Cells.Item(8, 17).Value = 723.7702819919538
$ws.Cells.Item(8, 18).Value = 6513.932537927583
$ws.Cells.Item(8, 19).Value = 0.105151951285593
$ws.Cells.Item(8, 20).Value = 0.105151951285593

$ws.Cells.Item(9, 9).Value = 0.4168441980730721
$ws.Cells.Item(9, 10).Value = 0.4168441980730722
$ws.Cells.Item(9, 13).Value = 0.4370123333333333
$ws.Cells.Item(9, 14).Value = 1.311037
$ws.Cells.Item(9, 15).Value = 0.005091641424883797
$ws.Cells.Item(9, 16).Value = 0.005091641424883797
$ws.Cells.Item(9, 17).Value = 14.60881478634356
$ws.Cells.Item(9, 18).Value = 131.479333077092
$ws.Cells.Item(9, 19).Value = 0.00212242118663132
$ws.Cells.Item(9, 20).Value = 0.002122421186631321

$ws.Cells.Item(10, 7).Value = 31.78201566666667
$ws.Cells.Item(10, 8).Value = 95.346047
$ws.Cells.Item(10, 9).Value = 0.3963089763847458
$ws.Cells.Item(10, 10).Value = 0.3963089763847459
$ws.Cells.Item(10, 13).Value = 4.114675333333333
$ws.Cells.Item(10, 14).Value = 12.344026
$ws.Cells.Item(10, 15).Value = 0.04794018332925969
$ws.Cells.Item(10, 16).Value = 0.04794018332925969
$ws.Cells.Item(10, 17).Value = 130.7726759072469
$ws.Cells.Item(10, 18).Value = 1176.954083165222
$ws.Cells.Item(10, 19).Value = 0.01899912498291596
$ws.Cells.Item(10, 20).Value = 0.01899912498291597

$ws.Cells.Item(11, 7).Value = 31.78201566666667
$ws.Cells.Item(11, 8).Value = 95.346047
$ws.Cells.Item(11, 9).Value = 0.3963089763847458
$ws.Cells.Item(11, 10).Value = 0.3963089763847459
$ws.Cells.Item(11, 13).Value = 59.62659933333333
$ws.Cells.Item(11, 15).Value = 0.6947109727426806
$ws.Cells.Item(11, 16).Value = 0.6947109727426806
$ws.Cells.Item(11, 17).Value = 1895.053514162056
$ws.Cells.Item(11, 18).Value = 17055.48162745851
$ws.Cells.Item(11, 19).Value = 0.2753201944909028
$ws.Cells.Item(11, 20).Value = 0.2753201944909029

$ws.Cells.Item(12, 7).Value = 31.78201566666667
$ws.Cells.Item(12, 8).Value = 95.346047
$ws.Cells.Item(12, 9).Value = 0.3963089763847458
$ws.Cells.Item(12, 10).Value = 0.3963089763847459
$ws.Cells.Item(12, 13).Value = 21.65107466666666
$ws.Cells.Item(12, 14).Value = 64.95322399999999
$ws.Cells.Item(12, 15).Value = 0.252257202503176
$ws.Cells.Item(12, 16).Value = 0.252257202503176
$ws.Cells.Item(12, 17).Value = 688.1147942561696
$ws.Cells.Item(12, 18).Value = 6193.033148305527
$ws.Cells.Item(12, 19).Value = 0.09997179370971322
$ws.Cells.Item(12, 20).Value = 0.09997179370971326

$ws.Cells.Item(13, 7).Value = 31.78201566666667
$ws.Cells.Item(13, 8).Value = 95.346047
$ws.Cells.Item(13, 9).Value = 0.3963089763847458
$ws.Cells.Item(13, 10).Value = 0.3963089763847459
$ws.Cells.Item(13, 13).Value = 0.4370123333333333
$ws.Cells.Item(13, 14).Value = 1.311037
$ws.Cells.Item(13, 15).Value = 0.005091641424883797
$ws.Cells.Item(13, 16).Value = 0.005091641424883797
$ws.Cells.Item(13, 17).Value = 13.88913282452656
$ws.Cells.Item(13, 18).Value = 125.002195420739
$ws.Cells.Item(13, 19).Value = 0.002017863201213866
$ws.Cells.Item(13, 20).Value = 0.002017863201213867

$ws.Cells.Item(14, 7).Value = 3.796148333333333
$ws.Cells.Item(14, 8).Value = 11.388445
$ws.Cells.Item(14, 9).Value = 0.04733644574235969
$ws.Cells.Item(14, 10).Value = 0.04733644574235969
$ws.Cells.Item(14, 13).Value = 4.114675333333333
$ws.Cells.Item(14, 14).Value = 12.344026
$ws.Cells.Item(14, 15).Value = 0.04794018332925969
$ws.Cells.Item(14, 16).Value = 0.04794018332925969
$ws.Cells.Item(14, 17).Value = 15.61991790884111
$ws.Cells.Item(14, 18).Value = 140.57926117957
$ws.Cells.Item(14, 19).Value = 0.002269317887044277
$ws.Cells.Item(14, 20).Value = 0.002269317887044278

$ws.Cells.Item(15, 7).Value = 3.796148333333333
$ws.Cells.Item(15, 8).Value = 11.388445
$ws.Cells.Item(15, 9).Value = 0.04733644574235969
$ws.Cells.Item(15, 10).Value = 0.04733644574235969
$ws.Cells.Item(15, 13).Value = 59.62659933333333
$ws.Cells.Item(15, 15).Value = 0.6947109727426806
$ws.Cells.Item(15, 16).Value = 0.6947109727426806
$ws.Cells.Item(15, 17).Value = 226.3514156815677
$ws.Cells.Item(15, 18).Value = 2037.16274113411
$ws.Cells.Item(15, 19).Value = 0.03288514826785582
$ws.Cells.Item(15, 20).Value = 0.03288514826785582

$ws.Cells.Item(16, 7).Value = 3.796148333333333
$ws.Cells.Item(16, 8).Value = 11.388445
$ws.Cells.Item(16, 9).Value = 0.04733644574235969
$ws.Cells.Item(16, 10).Value = 0.04733644574235969
$ws.Cells.Item(16, 13).Value = 21.65107466666666
$ws.Cells.Item(16, 14).Value = 64.95322399999999
$ws.Cells.Item(16, 15).Value = 0.252257202503176
$ws.Cells.Item(16, 16).Value = 0.252257202503176
$ws.Cells.Item(16, 17).Value = 82.1906910107422
$ws.Cells.Item(16, 18).Value = 739.7162190966799
$ws.Cells.Item(16, 19).Value = 0.01194095937941103
$ws.Cells.Item(16, 20).Value = 0.01194095937941103

$ws.Cells.Item(17, 7).Value = 3.796148333333333
$ws.Cells.Item(17, 8).Value = 11.388445
$ws.Cells.Item(17, 9).Value = 0.04733644574235969
$ws.Cells.Item(17, 10).Value = 0.04733644574235969
$ws.Cells.Item(17, 13).Value = 0.4370123333333333
$ws.Cells.Item(17, 14).Value = 1.311037
$ws.Cells.Item(17, 15).Value = 0.005091641424883797
$ws.Cells.Item(17, 16).Value = 0.005091641424883797
$ws.Cells.Item(17, 17).Value = 1.658963640829444
$ws.Cells.Item(17, 18).Value = 14.930672767465
$ws.Cells.Item(17, 19).Value = 0.0002410202080485628
$ws.Cells.Item(17, 20).Value = 0.0002410202080485628
